# Fixing issues with experiments validation
# The "atac-seq" sheet gains three new columns (Experiment Alias, Project,
# Secondary Project) right after the "Sample Descriptor" column, pushing
# every existing column from B..Y to E..AB.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("atac-seq")

# Insert three new columns before column B; this shifts the existing
# B:Y columns (and their data/styles) to E:AB automatically.
$ws.Columns("B:D").Insert()

# Populate the headers of the three newly-inserted columns.
$ws.Range("B1").Value = "Experiment Alias"
$ws.Range("C1").Value = "Project"
$ws.Range("D1").Value = "Secondary Project"
